$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# Add the new translation entry as the next row after the last used row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "empty_watchlist"
$ws.Cells.Item($newRow, 2).Value = "Merkliste leeren"

# Move the active selection to just past the newly added row, matching the
# updated "range control" scaler position after the criteria changed.
$ws.Range("A" + ($newRow + 1)).Select()
